# Roll back the "hot fix" that had prefixed several tab names with
# "Project - ". Restore the shorter tab names, and restore the
# previously-active sheet/selection (the 4th tab, "Funder", at G19)
# instead of the first tab ("Project").

$wb = $excel.ActiveWorkbook

# Rename the sheets back to their short names.
$wb.Worksheets.Item("Project - Contact").Name = "Contact"
$wb.Worksheets.Item("Project - Publications").Name = "Publications"
$wb.Worksheets.Item("Project - Funder").Name = "Funder"

# Re-select the cell that was active on the "Funder" sheet and make
# that sheet the active tab (this also clears tabSelected from
# whichever sheet currently has it, i.e. "Project").
$ws4 = $wb.Worksheets.Item("Funder")
$ws4.Activate()
$ws4.Range("G19").Select()
